# Apply updates to column G ("Recorded By") on the active worksheet.
# The values are re-ordered comma-separated lists of recorder identities.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Recorded By" value becomes "System, backup@backdoor.com, system"
$rowsA = @(2, 28, 54)
foreach ($r in $rowsA) {
    $ws.Range("G$r").Value = "System, backup@backdoor.com, system"
}

# Rows whose "Recorded By" value becomes "System, dnasr281@gmail.com"
$rowsB = @(3, 6, 10, 12, 13, 14, 15, 18, 19, 20, 21, 22, 24, 26, 29, 32, 36, 38, 39, 40, 41, 44, 45, 46, 47, 48, 50, 52, 55, 58, 62, 64, 65, 66, 67, 70, 71, 72, 73, 74, 76, 78, 83, 84, 85, 86, 90, 92, 99, 101, 109, 110, 111, 112, 116, 118, 125, 127, 135, 136, 137, 138, 142, 144, 151, 153)
foreach ($r in $rowsB) {
    $ws.Range("G$r").Value = "System, dnasr281@gmail.com"
}
